# Add 2023 ("genders in batumi" / Batumi region trade data) as new column S,
# mirroring the formatting of the existing column R, and re-center the
# merged title row across the newly widened table (A1:S1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Bring column R's formatting over to the new column S.
# ---------------------------------------------------------------------

# Row 1 (title row) - S1 should look like the other title cells (Q1), and
# R1 should stop using the special "last column" style and use the same
# style as the rest of the title row too.
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("Q1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows (3-14) - copy R's number formatting down into S.
$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Fill in the new 2023 figures.
# ---------------------------------------------------------------------
$ws.Range("S3").Value = 2023

$ws.Range("S4").Value = 1124.4000000000001
$ws.Range("S5").Value = 170.8
$ws.Range("S6").Value = 7146
$ws.Range("S7").Value = 4928
$ws.Range("S8").Value = 650.20000000000005
$ws.Range("S9").Value = 35
$ws.Range("S10").Value = 38.9
$ws.Range("S11").Value = 135.69999999999999
$ws.Range("S12").Value = 10.7
$ws.Range("S13").Value = 1014
$ws.Range("S14").Value = 981.2

# ---------------------------------------------------------------------
# 3. Merge the (now wider) title row and center it.
# ---------------------------------------------------------------------
$titleRng = $ws.Range("A1:S1")
$titleRng.Merge() | Out-Null
$titleRng.IndentLevel = 0
$titleRng.HorizontalAlignment = -4108   # xlCenter
$titleRng.VerticalAlignment = -4108     # xlCenter

# ---------------------------------------------------------------------
# 4. Match the selection the author left behind on the new last column.
# ---------------------------------------------------------------------
$ws.Range("S3:S14").Select() | Out-Null

Write-Output "done"
